# Natmi following Dr Hou advice
# Updates Ligand-expressing cells / Receptor-expressing cells counts (1 -> 3)
# and recomputed downstream expression / specificity / weight metrics
# for rows 2-17 of the LR-pairs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change, in order, for each data row: E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$columns = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

$rowData = @(
    "2|3,28.95628266666667,86.868848,0.5491054194301004,0.5491054194301005,3,218.9722396666667,656.9167190000001,0.9460661569256378,0.9460661569256379,6340.62206794108,57065.59861146972,0.5194900539072755,0.5194900539072757",
    "3|3,28.95628266666667,86.868848,0.5491054194301004,0.5491054194301005,3,5.232695,15.698085,0.02260777739018391,0.02260777739018392,151.5193955284533,1363.67455975608,0.01241405308621928,0.01241405308621928",
    "4|3,28.95628266666667,86.868848,0.5491054194301004,0.5491054194301005,3,3.718869666666667,11.156609,0.01606731857429249,0.0160673185742925,107.6846412684925,969.1617714164321,0.008822651704853921,0.008822651704853926",
    "5|3,28.95628266666667,86.868848,0.5491054194301004,0.5491054194301005,3,3.531721333333333,10.595164,0.01525874710988573,0.01525874710988573,102.2655212278969,920.3896910510721,0.008378660731751636,0.00837866073175164",
    "6|3,12.691493,38.074479,0.2406720388519202,0.2406720388519202,3,218.9722396666667,656.9167190000001,0.9460661569256378,0.9460661569256379,2779.084646923823,25011.7618223144,0.2276916708760939,0.2276916708760939",
    "7|3,12.691493,38.074479,0.2406720388519202,0.2406720388519202,3,5.232695,15.698085,0.02260777739018391,0.02260777739018392,66.41071196363499,597.6964076727149,0.005441059878405905,0.005441059878405906",
    "8|3,12.691493,38.074479,0.2406720388519202,0.2406720388519202,3,3.718869666666667,11.156609,0.01606731857429249,0.0160673185742925,47.19800834241234,424.782075081711,0.003866954320158301,0.003866954320158302",
    "9|3,12.691493,38.074479,0.2406720388519202,0.2406720388519202,3,3.531721333333333,10.595164,0.01525874710988573,0.01525874710988573,44.82281657995067,403.405349219556,0.003672353777262044,0.003672353777262044",
    "10|3,0.4888703333333334,1.466611,0.009270573592685367,0.009270573592685367,3,218.9722396666667,656.9167190000001,0.9460661569256378,0.9460661569256379,107.0490317965899,963.4412861693091,0.008770575931328147,0.008770575931328149",
    "11|3,0.4888703333333334,1.466611,0.009270573592685367,0.009270573592685367,3,5.232695,15.698085,0.02260777739018391,0.02260777739018392,2.558109348881667,23.022984139935,0.0002095870640627483,0.0002095870640627483",
    "12|3,0.4888703333333334,1.466611,0.009270573592685367,0.009270573592685367,3,3.718869666666667,11.156609,0.01606731857429249,0.0160673185742925,1.818045053566556,16.362405482099,0.0001489532592800991,0.0001489532592800991",
    "13|3,0.4888703333333334,1.466611,0.009270573592685367,0.009270573592685367,3,3.531721333333333,10.595164,0.01525874710988573,0.01525874710988573,1.726553785467111,15.538984069204,0.0001414573380143708,0.0001414573380143708",
    "14|3,10.59691233333333,31.790737,0.200951968125294,0.200951968125294,3,218.9722396666667,656.9167190000001,0.9460661569256378,0.9460661569256379,2320.429627181323,20883.8666446319,0.1901138562109402,0.1901138562109402",
    "15|3,10.59691233333333,31.790737,0.200951968125294,0.200951968125294,3,5.232695,15.698085,0.02260777739018391,0.02260777739018392,55.45041018207166,499.053691638645,0.00454307736149598,0.004543077361495981",
    "16|3,10.59691233333333,31.790737,0.200951968125294,0.200951968125294,3,3.718869666666667,11.156609,0.01606731857429249,0.0160673185742925,39.40853583675923,354.676822530833,0.00322875929000017,0.00322875929000017",
    "17|3,10.59691233333333,31.790737,0.200951968125294,0.200951968125294,3,3.531721333333333,10.595164,0.01525874710988573,0.01525874710988573,37.42534135509645,336.828072195868,0.00306627526285768,0.00306627526285768"
)

foreach ($entry in $rowData) {
    $parts = $entry.Split("|")
    $rowNum = $parts[0]
    $values = $parts[1].Split(",")

    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$rowNum"
        $ws.Range($cellRef).Value = [double]$values[$i]
    }
}

Write-Host "Updated rows 2-17 with new ligand/receptor expressing cell counts and recomputed metrics"
